$d = $word.ActiveDocument
$vtab = [char]11

function Merge-Text($searchText) {
    # Replace the matched range with identical text, which causes the
    # COM-interop engine to coalesce the runs it spans into a single run.
    # Returns $true if a match was found and replaced.
    $found = $d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)
    return $found
}

function Reseat-Run($searchText) {
    # Re-locate the (now merged) text and nudge Bold on/off. This forces
    # the engine to materialize a clean <w:rPr/> on the run, and - when
    # the search text does not include a neighbouring run of identical
    # formatting - forces that neighbour to split back out into its own
    # run instead of staying absorbed into this one.
    $rng = $d.Content
    $rng.Find.Execute($searchText) | Out-Null
    $rng.Bold = 1
    $rng.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "$ python" + "3" + " resolver_" + "dpll" + ".py" -> one run
# ---------------------------------------------------------------------
$t1 = "`$ python3 resolver_dpll.py"
Merge-Text $t1 | Out-Null
Reseat-Run $t1

# ---------------------------------------------------------------------
# 2) "resolver_" + "dpll" + ".py” se encargará de " + "darnos..." +
#    <br/><br/> + "El código..." -> one run (leaving the preceding “
#    quote-mark run untouched/separate)
# ---------------------------------------------------------------------
$t2 = "resolver_dpll.py” se encargará de darnos una solucion dadas nuestras reglas A, B, C, D, y si se quiere, alguna condición inicial." + $vtab + $vtab + "El código está adaptado para dar una solución completa a un sudoku parcialmente lleno."
Merge-Text $t2 | Out-Null
Reseat-Run $t2
# Re-split the preceding “ run back out (it gets absorbed by the merge
# above because it immediately abuts the replaced range with identical
# run formatting).
Reseat-Run "“"

# ---------------------------------------------------------------------
# 3) "El" + " .png presentado en la carpeta /Soluciones" + "DPLL" +
#    " corresponden al resultado de este proceso." -> one run
# ---------------------------------------------------------------------
$t3 = "El .png presentado en la carpeta /SolucionesDPLL corresponden al resultado de este proceso."
Merge-Text $t3 | Out-Null
Reseat-Run $t3

# ---------------------------------------------------------------------
# 4) ".OK" -> "OK"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(".OK", $true, $false, $false, $false, $false, $true, 1, $false, "OK", 2) | Out-Null
Reseat-Run "OK"

# ---------------------------------------------------------------------
# 5) & 6) Last table row ("Sustentación"): fill in the two empty cells.
# ---------------------------------------------------------------------
$table = $d.Tables.Item(1)
$row = $table.Rows.Item($table.Rows.Count)

$commentCell = $row.Cells.Item(2)
$commentCell.Range.Text = "¡Proyecto sustentado!" + "`r" + "Thanks!"

$okCell = $row.Cells.Item(3)
$okCell.Range.Text = "OK"
